$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 125
$ws.Range("I55").Value = 95.09090999999999
$ws.Range("J55").Value = 289.5
$ws.Range("K55").Value = 95.09090999999999
$ws.Range("L55").Value = 289.5
$ws.Range("M55").Value = 118.90909
$ws.Range("N55").Value = -717.5
$ws.Range("H100").Value = 2551.3333
$ws.Range("J100").Value = 1994.75
$ws.Range("L100").Value = 1994.75
$ws.Range("N100").Value = -3076.75
$ws.Range("H115").Value = 1445
$ws.Range("I115").Value = 1445
$ws.Range("K115").Value = 4335
$ws.Range("M115").Value = -2768
$ws.Range("H132").Value = 1802.1428
$ws.Range("I132").Value = 1702.7059
$ws.Range("J132").Value = 2224.75
$ws.Range("K132").Value = 5108.1177
$ws.Range("L132").Value = 6674.25
$ws.Range("M132").Value = -2578.1177
$ws.Range("N132").Value = -11734.25
$ws.Range("H138").Value = 11276.837
$ws.Range("I138").Value = 9831.666999999999
$ws.Range("J138").Value = 11385.225
$ws.Range("K138").Value = 29495.001
$ws.Range("L138").Value = 34155.675
$ws.Range("M138").Value = -24355.001
$ws.Range("N138").Value = -44435.675

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9291.325000000001
$ws.Range("I32").Value = 8464.632
$ws.Range("J32").Value = 24998.5
$ws.Range("K32").Value = 8464.632
$ws.Range("L32").Value = 24998.5
$ws.Range("M32").Value = -8177.632
$ws.Range("N32").Value = -25572.5
$ws.Range("H61").Value = 3921.889
$ws.Range("I61").Value = 2059.8
$ws.Range("J61").Value = 6249.5
$ws.Range("K61").Value = 2059.8
$ws.Range("L61").Value = 6249.5
$ws.Range("M61").Value = -1847.8
$ws.Range("N61").Value = -6673.5
$ws.Range("H63").Value = 2499.6667
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 2499.6667
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H74").Value = 3005.9285
$ws.Range("I74").Value = 3005.9285
$ws.Range("K74").Value = 3005.9285
$ws.Range("M74").Value = -2131.9285
$ws.Range("H77").Value = 3005.9285
$ws.Range("I77").Value = 3005.9285
$ws.Range("K77").Value = 15029.6425
$ws.Range("M77").Value = -10661.6425
$ws.Range("H94").Value = 18000
$ws.Range("J94").Value = 18000
$ws.Range("L94").Value = 18000
$ws.Range("N94").Value = -19802
$ws.Range("H136").Value = 3921.889
$ws.Range("I136").Value = 2059.8
$ws.Range("J136").Value = 6249.5
$ws.Range("K136").Value = 6179.400000000001
$ws.Range("L136").Value = 18748.5
$ws.Range("M136").Value = -3629.400000000001
$ws.Range("N136").Value = -23848.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1340.1904
$ws.Range("I94").Value = 1213.1578
$ws.Range("J94").Value = 2547
$ws.Range("K94").Value = 1213.1578
$ws.Range("L94").Value = 2547
$ws.Range("M94").Value = -762.1578
$ws.Range("N94").Value = -3449
$ws.Range("H134").Value = 4500
$ws.Range("I134").Value = 4500
$ws.Range("K134").Value = 13500
$ws.Range("M134").Value = -10965

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 90000
$ws.Range("J9").Value = 90000
$ws.Range("L9").Value = 90000
$ws.Range("N9").Value = -90336
$ws.Range("H16").Value = 2238
$ws.Range("J16").Value = 4325.6665
$ws.Range("L16").Value = 4325.6665
$ws.Range("N16").Value = -4899.6665
$ws.Range("H31").Value = 3060.3333
$ws.Range("I31").Value = 2957.6428
$ws.Range("K31").Value = 2957.6428
$ws.Range("M31").Value = -2662.6428
$ws.Range("H34").Value = 3060.3333
$ws.Range("I34").Value = 2957.6428
$ws.Range("K34").Value = 2957.6428
$ws.Range("M34").Value = -2755.6428
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = ""
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("H113").Value = 2238
$ws.Range("J113").Value = 4325.6665
$ws.Range("L113").Value = 4325.6665
$ws.Range("N113").Value = -8665.666499999999

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 664.5
$ws.Range("I116").Value = 664.5
$ws.Range("K116").Value = 1993.5
$ws.Range("M116").Value = 1448.5
$ws.Range("H117").Value = 122.5
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = ""
$ws.Range("H128").Value = 499997
$ws.Range("I128").Value = 499997
$ws.Range("K128").Value = 1499991
$ws.Range("M128").Value = -1495011
$ws.Range("H129").Value = 2261.2
$ws.Range("I129").Value = 270
$ws.Range("J129").Value = 3588.6667
$ws.Range("K129").Value = 810
$ws.Range("L129").Value = 10766.0001
$ws.Range("M129").Value = 4190
$ws.Range("N129").Value = -20766.0001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7239
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H102").Value = 1135.6666
$ws.Range("I102").Value = 1135.6666
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1135.6666
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 486.3334
$ws.Range("N102").Value = ""
$ws.Range("H132").Value = 2009.909
$ws.Range("I132").Value = 1900.8572
$ws.Range("K132").Value = 5702.571599999999
$ws.Range("M132").Value = -3172.571599999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 40000
$ws.Range("J41").Value = 40000
$ws.Range("L41").Value = 40000
$ws.Range("N41").Value = -40876
$ws.Range("H92").Value = 23055.5
$ws.Range("I92").Value = 22111
$ws.Range("K92").Value = 22111
$ws.Range("M92").Value = -19615
$ws.Range("H132").Value = 3757.8572
$ws.Range("I132").Value = 2660
$ws.Range("J132").Value = 6502.5
$ws.Range("K132").Value = 7980
$ws.Range("L132").Value = 19507.5
$ws.Range("M132").Value = -5450
$ws.Range("N132").Value = -24567.5
$ws.Range("H136").Value = 3727
$ws.Range("I136").Value = 2856.7144
$ws.Range("J136").Value = 5250
$ws.Range("K136").Value = 8570.143199999999
$ws.Range("L136").Value = 15750
$ws.Range("M136").Value = -6020.143199999999
$ws.Range("N136").Value = -20850

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 625
$ws.Range("J2").Value = 1034.3334
$ws.Range("L2").Value = 1034.3334
$ws.Range("N2").Value = -1258.3334
$ws.Range("H4").Value = 5050.5
$ws.Range("I4").Value = 7000
$ws.Range("J4").Value = 4400.6665
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 4400.6665
$ws.Range("M4").Value = -6887
$ws.Range("N4").Value = -4626.6665
$ws.Range("H132").Value = 3236.4546
$ws.Range("I132").Value = 3300.7778
$ws.Range("J132").Value = 2947
$ws.Range("K132").Value = 9902.3334
$ws.Range("L132").Value = 8841
$ws.Range("M132").Value = -7372.3334
$ws.Range("N132").Value = -13901
